# Edit testResults.xlsx per commit: rename legacy/newer weights rows to
# include the model-size suffix, and append the testSamples30-7.py group-
# voting + single-image evaluation rows (156-162).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits: rename existing entries to include model variant ---
$ws.Range("A139").Value = "legacy weights, small2d"
$ws.Range("A141").Value = "newer weights, small2d"
$ws.Range("A143").Value = "legacy weights, large2d"
$ws.Range("A145").Value = "newer weights, large2d"

# --- New rows 156-162: testSamples30-7.py changes (group voting + single image) ---
$ws.Range("A156").Value = "testSamples30-7.py -- Changed f1 and recall to have multilabel scores in the std respectively. ALSO made grouped2D group patient slices togheter rather than treat each slice individually"
$ws.Range("A157").Value = "viewing performance on grouped2d majourity voting"
$ws.Range("A158").Value = "Tests/0--/foldn5"
$ws.Range("B158").Value = 100
$ws.Range("C158").Value = 8
$ws.Range("D158").Value = 0.001
$ws.Range("E158").Value = 0.2
$ws.Range("F158").Value = 0.01
$ws.Range("G158").Value = "python testSamples30-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine='viewing performance on grouped2d majourity voting' -hasBackground=f -usesLargestBox=f -segmentsMultiple=13 -dropoutRate=0.2 -grouped2D=t -modelChosen='Small2D'"
$ws.Range("I158").Value = "[None, None]"
$ws.Range("L158").Value = "[None, None]"
$ws.Range("O158").Value = "{0: 5.2, 1: 8.0, 2: 3.8}"
$ws.Range("P158").Value = 0.34019607843137262
$ws.Range("Q158").Value = 0.284400871459695
$ws.Range("R158").Value = 0.30613756613756621
$ws.Range("S158").Value = 0.14579888192247431
$ws.Range("T158").Value = "{0: 0.2, 1: 0.4620915032679738, 2: 0.1911111111111111}"
$ws.Range("U158").Value = "{0: 0.22000000000000003, 1: 0.5095238095238096, 2: 0.18888888888888888}"
$ws.Range("V158").Value = "[0.35294117647058826, 0.5294117647058824, 0.23529411764705882, 0.3333333333333333, 0.25]"
$ws.Range("W158").Value = "[[0.0, 0.588235294117647, 0.2222222222222222], [0.5714285714285715, 0.5, 0.5333333333333333], [0.0, 0.4444444444444444, 0.0], [0.42857142857142855, 0.3333333333333333, 0.19999999999999998], [0.0, 0.4444444444444444, 0.0]]"
$ws.Range("X158").Value = "[[0.0, 0.7142857142857143, 0.16666666666666666], [0.5, 0.42857142857142855, 0.6666666666666666], [0.0, 0.5714285714285714, 0.0], [0.6, 0.5, 0.1111111111111111], [0.0, 0.3333333333333333, 0.0]]"
$ws.Range("Y158").Value = "[17, 16, 13, 14, 28]"
$ws.Range("A159").Value = "viewing performance on grouped2d average voting"
$ws.Range("A160").Value = "Tests/0--/foldn5"
$ws.Range("B160").Value = 100
$ws.Range("C160").Value = 8
$ws.Range("D160").Value = 0.001
$ws.Range("E160").Value = 0.2
$ws.Range("F160").Value = 0.01
$ws.Range("G160").Value = "python testSamples30-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine='viewing performance on grouped2d average voting' -hasBackground=f -usesLargestBox=f -segmentsMultiple=13 -dropoutRate=0.2 -grouped2D=t -modelChosen='Small2D'"
$ws.Range("I160").Value = "[None, None]"
$ws.Range("L160").Value = "[None, None]"
$ws.Range("O160").Value = "{0: 1.4, 1: 11.0, 2: 4.6}"
$ws.Range("P160").Value = 0.44828431372549021
$ws.Range("Q160").Value = 0.35934922777028039
$ws.Range("R160").Value = 0.422010582010582
$ws.Range("S160").Value = 0.1157272360436034
$ws.Range("T160").Value = "{0: 0.17142857142857143, 1: 0.5678889531521111, 2: 0.3387301587301587}"
$ws.Range("U160").Value = "{0: 0.13999999999999999, 1: 0.7571428571428571, 2: 0.36888888888888893}"
$ws.Range("V160").Value = "[0.4117647058823529, 0.5882352941176471, 0.47058823529411764, 0.3333333333333333, 0.4375]"
$ws.Range("W160").Value = "[[0.0, 0.631578947368421, 0.2222222222222222], [0.5714285714285715, 0.6153846153846153, 0.5714285714285715], [0.0, 0.5714285714285714, 0.5], [0.28571428571428575, 0.4210526315789474, 0.19999999999999998], [0.0, 0.6, 0.2]]"
$ws.Range("X160").Value = "[[0.0, 0.8571428571428571, 0.16666666666666666], [0.5, 0.5714285714285714, 0.6666666666666666], [0.0, 0.8571428571428571, 0.4], [0.2, 1.0, 0.1111111111111111], [0.0, 0.5, 0.5]]"
$ws.Range("Y160").Value = "[17, 16, 13, 14, 28]"
$ws.Range("A161").Value = "viewing performance on single image"
$ws.Range("A162").Value = "Tests/0--/foldn5"
$ws.Range("B162").Value = 100
$ws.Range("C162").Value = 8
$ws.Range("D162").Value = 0.001
$ws.Range("E162").Value = 0.2
$ws.Range("F162").Value = 0.01
$ws.Range("G162").Value = "python testSamples30-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine='viewing performance on single image' -hasBackground=f -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f -modelChosen='Small2D'"
$ws.Range("I162").Value = "[None, None]"
$ws.Range("L162").Value = "[None, None]"
$ws.Range("O162").Value = "{0: 4.2, 1: 9.8, 2: 3.0}"
$ws.Range("P162").Value = 0.49411764705882361
$ws.Range("Q162").Value = 0.43601589443694722
$ws.Range("R162").Value = 0.46341269841269839
$ws.Range("S162").Value = 0.13268203919406141
$ws.Range("T162").Value = "{0: 0.37222222222222223, 1: 0.6069365721997302, 2: 0.3288888888888889}"
$ws.Range("U162").Value = "{0: 0.4, 1: 0.7035714285714285, 2: 0.2866666666666667}"
$ws.Range("V162").Value = "[0.35294117647058826, 0.6470588235294118, 0.5294117647058824, 0.47058823529411764, 0.47058823529411764]"
$ws.Range("W162").Value = "[[0.0, 0.4210526315789474, 0.4444444444444445], [0.75, 0.7000000000000001, 0.33333333333333337], [0.22222222222222224, 0.6153846153846153, 0.6666666666666666], [0.4444444444444445, 0.631578947368421, 0.0], [0.4444444444444445, 0.6666666666666666, 0.2]]"
$ws.Range("X162").Value = "[[0.0, 0.5, 0.4], [0.75, 0.875, 0.2], [0.25, 0.5714285714285714, 0.6666666666666666], [0.5, 0.8571428571428571, 0.0], [0.5, 0.7142857142857143, 0.16666666666666666]]"
$ws.Range("Y162").Value = "[19, 25, 35, 17, 21]"

# Restore the view roughly where the author left it: scrolled down to the
# new rows, with the cursor parked just past the last data row.
$w = $excel.ActiveWindow
$w.ScrollRow = 141
$w.ScrollColumn = 1
$ws.Range("N167").Select()
